# Balance.xlsx update: add "Operative max forward Xcg MAC" row on the
# GLOBAL RESULTS sheet, and reorder the SFORZA / TORENBEEK_1982 comparison
# rows (SFORZA now listed before TORENBEEK_1982) on the FUSELAGE and WING
# sheets.

$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS: insert "Operative max forward Xcg MAC" row ---------
$wsGlobal = $wb.Worksheets.Item("GLOBAL RESULTS")

# Row 63 currently holds "Max aft Xcg MAC"; push it down to make room for
# the new "Operative max forward Xcg MAC" row right after "Max forward Xcg
# MAC" (row 62).
$wsGlobal.Rows.Item(63).Insert()

$wsGlobal.Range("A63").Value = "Operative max forward Xcg MAC"
$wsGlobal.Range("B63").Value = "%"
$wsGlobal.Range("C63").Value = 25.986504960387446

# --- FUSELAGE: swap SFORZA / TORENBEEK_1982 rows -------------------------
$wsFuselage = $wb.Worksheets.Item("FUSELAGE")

$wsFuselage.Range("A11").Value = "SFORZA"
$wsFuselage.Range("C11").Value = 17.143322222222217
$wsFuselage.Range("A12").Value = "TORENBEEK_1982"
$wsFuselage.Range("C12").Value = 16.8345

# --- WING: swap SFORZA / TORENBEEK_1982 rows (two comparison blocks) ----
$wsWing = $wb.Worksheets.Item("WING")

$wsWing.Range("A11").Value = "SFORZA"
$wsWing.Range("C11").Value = 4.3629715646212155
$wsWing.Range("A12").Value = "TORENBEEK_1982"
$wsWing.Range("C12").Value = 3.5180298935880643

$wsWing.Range("A15").Value = "SFORZA"
$wsWing.Range("C15").Value = 4.998846772296348
$wsWing.Range("A16").Value = "TORENBEEK_1982"
$wsWing.Range("C16").Value = 6.114221148470394
